$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Asia)
$ws.Range("B3").Value = 217433783
$ws.Range("C3").Value = 3162
$ws.Range("D3").Value = 1546556
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = 201147206
$ws.Range("G3").Value = 26811
$ws.Range("H3").Value = 14740021
$ws.Range("I3").Value = 15346

# Row 4 (Europe)
$ws.Range("B4").Value = 249364590
$ws.Range("C4").Value = 189
$ws.Range("D4").Value = 2061026
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 245341950
$ws.Range("G4").Value = 5581
$ws.Range("H4").Value = 1961614

# Row 6 (Australia/Oceania)
$ws.Range("F6").Value = 14188518
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 148580
